$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 229
$ws.Range("I2").Value = 211.25
$ws.Range("K2").Value = 211.25
$ws.Range("M2").Value = -98.25

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 616.6667
$ws.Range("I18").Value = 616.6667
$ws.Range("K18").Value = 616.6667
$ws.Range("M18").Value = -332.6667

# Row 29 (Leve Item ID 4575)
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 9
$ws.Range("I38").Value = 9
$ws.Range("K38").Value = 27
$ws.Range("M38").Value = 345

# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 200
$ws.Range("I58").Value = 200
$ws.Range("K58").Value = 600
$ws.Range("M58").Value = -450

# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 6000
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -5752
$ws.Range("N64").ClearContents()

# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 6000
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -5142
$ws.Range("N67").ClearContents()

# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 473.9091
$ws.Range("I92").Value = 517
$ws.Range("J92").Value = 43
$ws.Range("K92").Value = 517
$ws.Range("L92").Value = 43
$ws.Range("M92").Value = 731
$ws.Range("N92").Value = -2539

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 3250
$ws.Range("K100").Value = 3250
$ws.Range("M100").Value = -2709

$ws = $wb.Worksheets.Item("ARM")
# Row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 49790
$ws.Range("J37").Value = 49790
$ws.Range("L37").Value = 49790
$ws.Range("N37").Value = -50336

# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2878.1
$ws.Range("I102").Value = 2878.1
$ws.Range("K102").Value = 2878.1
$ws.Range("M102").Value = -1256.1

# Row 113 (Leve Item ID 26002)
$ws.Range("H113").Value = 39990
$ws.Range("J113").Value = 39990
$ws.Range("L113").Value = 39990
$ws.Range("N113").Value = -48668

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 4115.222
$ws.Range("J122").Value = 4013.8
$ws.Range("L122").Value = 12041.4
$ws.Range("N122").Value = -16941.4

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 3028.4285
$ws.Range("I99").Value = 3375
$ws.Range("J99").Value = 2566.3333
$ws.Range("K99").Value = 3375
$ws.Range("L99").Value = 2566.3333
$ws.Range("M99").Value = -1877
$ws.Range("N99").Value = -5562.3333

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4736.75
$ws.Range("I134").Value = 4736.75
$ws.Range("K134").Value = 14210.25
$ws.Range("M134").Value = -11675.25

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 453.2
$ws.Range("I22").Value = 369.07693
$ws.Range("K22").Value = 369.07693
$ws.Range("M22").Value = -19.07693

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2276.125
$ws.Range("I31").Value = 2227.8667
$ws.Range("K31").Value = 2227.8667
$ws.Range("M31").Value = -1932.8667

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2276.125
$ws.Range("I34").Value = 2227.8667
$ws.Range("K34").Value = 2227.8667
$ws.Range("M34").Value = -2025.8667

# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 13959.8
$ws.Range("I62").Value = 3599.6667
$ws.Range("K62").Value = 3599.6667
$ws.Range("M62").Value = -2975.6667

# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 13959.8
$ws.Range("I65").Value = 3599.6667
$ws.Range("K65").Value = 17998.3335
$ws.Range("M65").Value = -14878.3335

# Row 88 (Leve Item ID 10608)
$ws.Range("H88").Value = 17332.834
$ws.Range("J88").Value = 17332.834
$ws.Range("L88").Value = 17332.834
$ws.Range("N88").Value = -18144.834

# Row 91 (Leve Item ID 10608)
$ws.Range("H91").Value = 17332.834
$ws.Range("J91").Value = 17332.834
$ws.Range("L91").Value = 17332.834
$ws.Range("N91").Value = -20140.834

$ws = $wb.Worksheets.Item("CUL")
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 60000
$ws.Range("N39").Value = -60588

# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 11500
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 11714.286
$ws.Range("K55").Value = 30000
$ws.Range("L55").Value = 35142.858
$ws.Range("M55").Value = -29823
$ws.Range("N55").Value = -35496.858

$ws = $wb.Worksheets.Item("GSM")
# Row 11 (Leve Item ID 4422)
$ws.Range("H11").Value = 7625000
$ws.Range("J11").Value = 500000
$ws.Range("L11").Value = 500000
$ws.Range("N11").Value = -500278

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1401.3334
$ws.Range("I102").Value = 1401.3334
$ws.Range("K102").Value = 1401.3334
$ws.Range("M102").Value = 220.6666

# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 57500
$ws.Range("I130").Value = 25000
$ws.Range("K130").Value = 25000
$ws.Range("M130").Value = -19980

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 3908.3333
$ws.Range("I7").Value = 3675
$ws.Range("J7").Value = 4375
$ws.Range("K7").Value = 3675
$ws.Range("L7").Value = 4375
$ws.Range("M7").Value = -3563
$ws.Range("N7").Value = -4599

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 4266
$ws.Range("I22").Value = 3477.5
$ws.Range("J22").Value = 4791.6665
$ws.Range("K22").Value = 3477.5
$ws.Range("L22").Value = 4791.6665
$ws.Range("M22").Value = -3182.5
$ws.Range("N22").Value = -5381.6665

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 4266
$ws.Range("I27").Value = 3477.5
$ws.Range("J27").Value = 4791.6665
$ws.Range("K27").Value = 3477.5
$ws.Range("L27").Value = 4791.6665
$ws.Range("M27").Value = -3370.5
$ws.Range("N27").Value = -5005.6665

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 3908.3333
$ws.Range("I126").Value = 3675
$ws.Range("J126").Value = 4375
$ws.Range("K126").Value = 11025
$ws.Range("L126").Value = 13125
$ws.Range("M126").Value = -8555
$ws.Range("N126").Value = -18065

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5999
$ws.Range("J132").Value = 5999
$ws.Range("L132").Value = 17997
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 6925
$ws.Range("J96").Value = 3700
$ws.Range("L96").Value = 3700
$ws.Range("N96").Value = -6446

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1998.5834
$ws.Range("I132").Value = 1293.3334
$ws.Range("K132").Value = 3880.0002
$ws.Range("M132").Value = -1350.0002

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1919.909
$ws.Range("I136").Value = 1811.9
$ws.Range("K136").Value = 5435.700000000001
$ws.Range("M136").Value = -2885.700000000001
